$wb = $excel.ActiveWorkbook

# Row 43 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7937165
$ws.Range("I43").Value = 466.33334
$ws.Range("K43").Value = 466.33334
$ws.Range("M43").Value = -397.33334

# Row 129 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 918.91895
$ws.Range("J129").Value = 923.6111
$ws.Range("L129").Value = 2770.8333
$ws.Range("N129").Value = -12770.8333

# Row 135 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 258.53845
$ws.Range("I135").Value = 263.6
$ws.Range("J135").Value = 255.375
$ws.Range("K135").Value = 2372.4
$ws.Range("L135").Value = 2298.375
$ws.Range("M135").Value = 162.5999999999999
$ws.Range("N135").Value = -7368.375

# Row 137 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1905.8462
$ws.Range("I137").Value = 1675.8422
$ws.Range("J137").Value = 2530.1428
$ws.Range("K137").Value = 5027.5266
$ws.Range("L137").Value = 7590.428400000001
$ws.Range("M137").Value = -2477.5266
$ws.Range("N137").Value = -12690.4284

# Row 138 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1116777.1
$ws.Range("J138").Value = 1404522.9
$ws.Range("L138").Value = 4213568.699999999
$ws.Range("N138").Value = -4223848.699999999

# Row 139 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 67326.664
$ws.Range("J139").Value = 67326.664
$ws.Range("L139").Value = 67326.664
$ws.Range("N139").Value = -77606.664

# Row 16 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 52 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()

# Row 61 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2401.8572
$ws.Range("I61").Value = 2199.75
$ws.Range("K61").Value = 2199.75
$ws.Range("M61").Value = -1987.75

# Row 74 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1250.4445
$ws.Range("I74").Value = 1330
$ws.Range("J74").Value = 614
$ws.Range("K74").Value = 1330
$ws.Range("L74").Value = 614
$ws.Range("M74").Value = -456
$ws.Range("N74").Value = -2362

# Row 77 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1250.4445
$ws.Range("I77").Value = 1330
$ws.Range("J77").Value = 614
$ws.Range("K77").Value = 6650
$ws.Range("L77").Value = 3070
$ws.Range("M77").Value = -2282
$ws.Range("N77").Value = -11806

# Row 136 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2401.8572
$ws.Range("I136").Value = 2199.75
$ws.Range("K136").Value = 6599.25
$ws.Range("M136").Value = -4049.25

# Row 37 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 493.2
$ws.Range("I37").Value = 493.2
$ws.Range("K37").Value = 493.2
$ws.Range("M37").Value = -356.2

# Row 57 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 24209
$ws.Range("I57").Value = 24209
$ws.Range("K57").Value = 24209
$ws.Range("M57").Value = -23489

# Row 136 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 24209
$ws.Range("I136").Value = 24209
$ws.Range("K136").Value = 24209
$ws.Range("M136").Value = -19109

# Row 31 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1055.7534
$ws.Range("I31").Value = 728.43396
$ws.Range("J31").Value = 1923.15
$ws.Range("K31").Value = 728.43396
$ws.Range("L31").Value = 1923.15
$ws.Range("M31").Value = -433.43396
$ws.Range("N31").Value = -2513.15

# Row 34 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1055.7534
$ws.Range("I34").Value = 728.43396
$ws.Range("J34").Value = 1923.15
$ws.Range("K34").Value = 728.43396
$ws.Range("L34").Value = 1923.15
$ws.Range("M34").Value = -526.43396
$ws.Range("N34").Value = -2327.15

# Row 58 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 933.4
$ws.Range("I58").Value = 933.4
$ws.Range("K58").Value = 933.4
$ws.Range("M58").Value = -730.4

# Row 86 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4182606.2
$ws.Range("I86").Value = 6064044
$ws.Range("J86").Value = 43443.2
$ws.Range("K86").Value = 6064044
$ws.Range("L86").Value = 43443.2
$ws.Range("M86").Value = -6062921
$ws.Range("N86").Value = -45689.2

# Row 89 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4182606.2
$ws.Range("I89").Value = 6064044
$ws.Range("J89").Value = 43443.2
$ws.Range("K89").Value = 30320220
$ws.Range("L89").Value = 217216
$ws.Range("M89").Value = -30314604
$ws.Range("N89").Value = -228448

# Row 107 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 540.4
$ws.Range("I107").Value = 169.83333
$ws.Range("J107").Value = 699.2143
$ws.Range("K107").Value = 169.83333
$ws.Range("L107").Value = 699.2143
$ws.Range("M107").Value = 1750.16667
$ws.Range("N107").Value = -4539.2143

# Row 136 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 933.4
$ws.Range("I136").Value = 933.4
$ws.Range("K136").Value = 2800.2
$ws.Range("M136").Value = -250.1999999999998

# Row 32 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2299.8
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2299.8
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 6899.400000000001
$ws.Range("N32").Value = -7465.400000000001
$ws.Range("M32").ClearContents()

# Row 113 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 648.6087
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 648.6087
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1945.8261
$ws.Range("N113").Value = -6285.8261
$ws.Range("M113").ClearContents()

# Row 131 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 38463290
$ws.Range("J131").Value = 2006.6818
$ws.Range("L131").Value = 6020.0454
$ws.Range("N131").Value = -16100.0454

# Row 136 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2374.7273
$ws.Range("I136").Value = 915
$ws.Range("J136").Value = 4126.4
$ws.Range("K136").Value = 2745
$ws.Range("L136").Value = 12379.2
$ws.Range("M136").Value = 2355
$ws.Range("N136").Value = -22579.2

# Row 137 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 13274
$ws.Range("I137").Value = 2618
$ws.Range("J137").Value = 20885.428
$ws.Range("K137").Value = 7854
$ws.Range("L137").Value = 62656.284
$ws.Range("M137").Value = -2754
$ws.Range("N137").Value = -72856.284

# Row 139 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2049.6667
$ws.Range("I139").Value = 2097.2942
$ws.Range("J139").Value = 1999.0625
$ws.Range("K139").Value = 6291.882599999999
$ws.Range("L139").Value = 5997.1875
$ws.Range("M139").Value = -1151.882599999999
$ws.Range("N139").Value = -16277.1875

# Row 141 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2261.25
$ws.Range("I141").Value = 2261.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6783.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1603.75
$ws.Range("N141").ClearContents()

# Row 97 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1896.8462
$ws.Range("I97").Value = 2059.9092
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 2059.9092
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -1563.9092
$ws.Range("N97").Value = -1992

# Row 100 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# Row 16 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1760.8182
$ws.Range("I16").Value = 1760.8182
$ws.Range("K16").Value = 1760.8182
$ws.Range("M16").Value = -1590.8182

# Row 40 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2726
$ws.Range("I40").Value = 2726
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2726
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2590
$ws.Range("N40").ClearContents()

# Row 122 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 18893616
$ws.Range("I122").Value = 25760476
$ws.Range("K122").Value = 77281428
$ws.Range("M122").Value = -77278978

# Row 136 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 13112.667
$ws.Range("I136").Value = 21347.6
$ws.Range("K136").Value = 64042.8
$ws.Range("M136").Value = -61492.8

# Row 141 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 50715
$ws.Range("J141").Value = 50715
$ws.Range("L141").Value = 50715
$ws.Range("N141").Value = -61075

# Row 62 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26323042
$ws.Range("I62").Value = 35719590
$ws.Range("J62").Value = 12700
$ws.Range("K62").Value = 35719590
$ws.Range("L62").Value = 12700
$ws.Range("M62").Value = -35718966
$ws.Range("N62").Value = -13948

# Row 65 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 26323042
$ws.Range("I65").Value = 35719590
$ws.Range("J65").Value = 12700
$ws.Range("K65").Value = 178597950
$ws.Range("L65").Value = 63500
$ws.Range("M65").Value = -178594830
$ws.Range("N65").Value = -69740

# Row 126 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2199
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 132 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4374.25
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
